$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.178152918815613
$ws.Range("B1").Value = 3.700886249542236
$ws.Range("C1").Value = 4.342037677764893
$ws.Range("D1").Value = 1.824571847915649
$ws.Range("E1").Value = 1.266490459442139
